# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback has been processed:
#   - The "Status" shown for every file changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview sheet + per-language sheets).
#   - Each per-language sheet (zh-cn, de-de) gets two new populated columns
#     for the files that were actually localized: "Latest Target File" (E)
#     and "Latest Handback File" (F), each rendered as a hyperlink just like
#     the existing "Source File Name" / "Latest Handoff File" columns.
#   - The "Latest Handback DateTime" column (G) is stamped with the time the
#     handback report was produced, for the rows that received a handback.

$wb = $excel.ActiveWorkbook

$oldStatusText = "Ready for handoff"
$newStatusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text everywhere it currently appears -------------
$wsOverview.Range("B2").Value = $newStatusText
$wsOverview.Range("C2").Value = $newStatusText
$wsOverview.Range("B3").Value = $newStatusText
$wsOverview.Range("C3").Value = $newStatusText

$wsZhCn.Range("B2").Value = $newStatusText
$wsZhCn.Range("B3").Value = $newStatusText

$wsDeDe.Range("B2").Value = $newStatusText
$wsDeDe.Range("B3").Value = $newStatusText

# --- 2. zh-cn sheet: add Latest Target File / Latest Handback File links ---
$zhCnHandoffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d66149f08223c0b6e816a66d2845a4aa26f08a0c/e2e/a.md"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dda4edf64615a451f5ecb7e668e869300c2c63f2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhCnXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), $zhCnHandoffUrl, "", "", "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhCnXlfUrl, "", "", $zhCnXlfName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), $zhCnHandoffUrl, "", "", "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhCnXlfUrl, "", "", $zhCnXlfName) | Out-Null

# Stamp the handback datetime for the rows that were just handed back.
$wsZhCn.Range("G2").Value = "2016-02-26 06:25:27"
$wsZhCn.Range("G3").Value = "2016-02-26 06:25:27"

# --- 3. de-de sheet: add Latest Target File / Latest Handback File links ---
$deDeHandoffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d66149f08223c0b6e816a66d2845a4aa26f08a0c/e2e/a.md"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58ec6bce44b5d0538fc2f19d2a1d8acf373e1bdc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deDeXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), $deDeHandoffUrl, "", "", "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deDeXlfUrl, "", "", $deDeXlfName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), $deDeHandoffUrl, "", "", "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deDeXlfUrl, "", "", $deDeXlfName) | Out-Null

# Stamp the handback datetime for the rows that were just handed back.
$wsDeDe.Range("G2").Value = "2016-02-26 06:25:49"
$wsDeDe.Range("G3").Value = "2016-02-26 06:25:49"
